# Weekly fruit/vegetable price update: a new observation (row) was recorded
# for Albahaca at "Vega Central Mapocho de Santiago" and inserted into the
# data table right after the existing row for 2021-01-29 (row 239), pushing
# every subsequent row down by one. The sheet's used range grows from
# A1:R346 to A1:R347.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 240 (existing rows 240.. shift to 241..)
$ws.Rows.Item(240).Insert()

# Populate the newly inserted row with the new observation's data
$ws.Cells.Item(240, 1).Value  = 9
$ws.Cells.Item(240, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(240, 3).Value  = "Metropolitana"
$ws.Cells.Item(240, 4).Value  = 44636
$ws.Cells.Item(240, 5).Value  = 13
$ws.Cells.Item(240, 6).Value  = 100112052
$ws.Cells.Item(240, 7).Value  = "Albahaca"
$ws.Cells.Item(240, 8).Value  = "Sin especificar"
$ws.Cells.Item(240, 9).Value  = "Primera"
$ws.Cells.Item(240, 10).Value = 400
$ws.Cells.Item(240, 11).Value = 3000
$ws.Cells.Item(240, 12).Value = 3500
$ws.Cells.Item(240, 13).Value = 3156
$ws.Cells.Item(240, 14).Value = "`$/docena de matas"
$ws.Cells.Item(240, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(240, 16).Value = 526
$ws.Cells.Item(240, 17).Value = 6
$ws.Cells.Item(240, 18).Value = "Hortaliza"
